$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values so they stay as text, matching original inlineStr cells
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "25.925.63"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "1.637.75"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "215.43"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "0.0638"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "19.61"
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").Value = "1.865.38"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").Value = "1.640.64"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "25.949.06"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "192.86"
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").Value = "6.29"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "1.79"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "144.17"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +3.82%  "
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("D30").Value = "1.24"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("E35").Value = "  +1.76%  "
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").Value = "1.137.24"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("D39").Value = "2.47"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "5.49"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "99.31"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.798"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.774.96"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0115"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "56.65"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.0531"
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.48"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.70"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.414"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.0963"
$ws.Range("E51").Value = "  -0.98%  "

# Clear the temporary text-number-format style so no stray style index remains on these cells
$ws.Range("D5").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

